$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.18"
$ws.Range("E2").Value = "'-2.30%"
$ws.Range("D3").Value = "'44.49"
$ws.Range("E3").Value = "'1.18%"
$ws.Range("D4").Value = "'5.605"
$ws.Range("E4").Value = "'-2.21%"
$ws.Range("D5").Value = "'0.08062"
$ws.Range("E5").Value = "'-3.33%"
$ws.Range("D6").Value = "'1.903"
$ws.Range("E6").Value = "'-3.17%"
$ws.Range("D7").Value = "'4.299"
$ws.Range("E7").Value = "'-4.84%"
$ws.Range("D8").Value = "'2.670"
$ws.Range("E8").Value = "'-7.28%"
$ws.Range("D9").Value = "'0.9453"
$ws.Range("E9").Value = "'0.22%"
$ws.Range("D10").Value = "'0.1160"
$ws.Range("E10").Value = "'-7.05%"
$ws.Range("D11").Value = "'0.1844"
$ws.Range("E11").Value = "'-6.95%"
$ws.Range("D12").Value = "'0.09851"
$ws.Range("E12").Value = "'-8.19%"
$ws.Range("E13").Value = "'-9.91%"
$ws.Range("D14").Value = "'0.1067"
$ws.Range("E14").Value = "'0.00%"
$ws.Range("D15").Value = "'0.001283"
$ws.Range("E15").Value = "'-1.19%"
$ws.Range("D16").Value = "'0.04217"
$ws.Range("E16").Value = "'-4.39%"
$ws.Range("D17").Value = "'0.005990"
$ws.Range("E17").Value = "'0.82%"
$ws.Range("D18").Value = "'3.609"
$ws.Range("E18").Value = "'3.13%"
$ws.Range("D19").Value = "'0.3500"
$ws.Range("E19").Value = "'-0.14%"
$ws.Range("D20").Value = "'8.429"
$ws.Range("E20").Value = "'-3.16%"
$ws.Range("D21").Value = "'0.1372"
$ws.Range("E21").Value = "'1.45%"
$ws.Range("D22").Value = "'0.2655"
$ws.Range("E22").Value = "'-1.30%"
$ws.Range("D23").Value = "'0.001248"
$ws.Range("E23").Value = "'-0.52%"
$ws.Range("D24").Value = "'0.004493"
$ws.Range("E24").Value = "'2.98%"
$ws.Range("D25").Value = "'0.0001261"
$ws.Range("E25").Value = "'-0.05%"
$ws.Range("D26").Value = "'0.0003997"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("D38").Value = "'0.02638"
$ws.Range("E38").Value = "'-6.60%"
$ws.Range("D39").Value = "'0.05469"
$ws.Range("E39").Value = "'-9.20%"
$ws.Range("D40").Value = "'0.007628"
$ws.Range("E40").Value = "'-3.82%"
$ws.Range("D41").Value = "'0.1399"
$ws.Range("E41").Value = "'-2.00%"
$ws.Range("D42").Value = "'0.007351"
$ws.Range("E42").Value = "'-17.96%"
$ws.Range("D43").Value = "'0.002049"
$ws.Range("E43").Value = "'-5.64%"
$ws.Range("D44").Value = "'0.008844"
$ws.Range("E44").Value = "'-12.90%"
$ws.Range("D45").Value = "'0.00006934"
$ws.Range("E45").Value = "'-0.98%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("D47").Value = "'0.003670"
$ws.Range("E47").Value = "'15.12%"
$ws.Range("D48").Value = "'0.002275"
$ws.Range("E48").Value = "'0.16%"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("E50").Value = "'0.06%"
